$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 100.888885
$ws.Range("I4").Value = 73.28570999999999
$ws.Range("J4").Value = 197.5
$ws.Range("K4").Value = 73.28570999999999
$ws.Range("L4").Value = 197.5
$ws.Range("M4").Value = 40.71429000000001
$ws.Range("N4").Value = -425.5

$ws.Range("H6").Value = 229.4
$ws.Range("I6").Value = 229.4
$ws.Range("K6").Value = 688.2
$ws.Range("M6").Value = -576.2

$ws.Range("H8").Value = 21
$ws.Range("I8").Value = 21
$ws.Range("K8").Value = 63
$ws.Range("M8").Value = 76

$ws.Range("H21").Value = 1001
$ws.Range("I21").Value = 1001
$ws.Range("K21").Value = 1001
$ws.Range("M21").Value = -533

$ws.Range("H23").Value = 1001
$ws.Range("I23").Value = 1001
$ws.Range("K23").Value = 1001
$ws.Range("M23").Value = -767

$ws.Range("H40").Value = 5600
$ws.Range("I40").Value = 5333.3335
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 5333.3335
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -5158.3335
$ws.Range("N40").Value = -6350

$ws.Range("H51").Value = 10168.462
$ws.Range("I51").Value = 8687
$ws.Range("J51").Value = 12538.8
$ws.Range("K51").Value = 8687
$ws.Range("L51").Value = 12538.8
$ws.Range("M51").Value = -8203
$ws.Range("N51").Value = -13506.8

$ws.Range("H138").Value = 11043.444
$ws.Range("I138").Value = 13798.333
$ws.Range("K138").Value = 41394.999
$ws.Range("M138").Value = -36254.999

$ws.Range("H141").Value = 2655
$ws.Range("I141").Value = 2497.6667
$ws.Range("K141").Value = 7493.000100000001
$ws.Range("M141").Value = -2313.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3420.9
$ws.Range("I2").Value = 3467.6667
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 3467.6667
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -3354.6667
$ws.Range("N2").Value = -3226

$ws.Range("H61").Value = 1224.6666
$ws.Range("I61").Value = 1201.375
$ws.Range("K61").Value = 1201.375
$ws.Range("M61").Value = -989.375

$ws.Range("H116").Value = 3420.9
$ws.Range("I116").Value = 3467.6667
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 3467.6667
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = -1173.6667
$ws.Range("N116").Value = -7588

$ws.Range("H132").Value = 2122.4
$ws.Range("I132").Value = 2037.3334
$ws.Range("K132").Value = 6112.0002
$ws.Range("M132").Value = -3582.0002

$ws.Range("H136").Value = 1224.6666
$ws.Range("I136").Value = 1201.375
$ws.Range("K136").Value = 3604.125
$ws.Range("M136").Value = -1054.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3420.9
$ws.Range("I3").Value = 3467.6667
$ws.Range("J3").Value = 3000
$ws.Range("K3").Value = 3467.6667
$ws.Range("L3").Value = 3000
$ws.Range("M3").Value = -3353.6667
$ws.Range("N3").Value = -3228

$ws.Range("H20").Value = 1054.1428
$ws.Range("I20").Value = 979.8333
$ws.Range("K20").Value = 979.8333
$ws.Range("M20").Value = -732.8333

$ws.Range("H22").Value = 279.5
$ws.Range("I22").Value = 276.57144
$ws.Range("K22").Value = 276.57144
$ws.Range("M22").Value = -103.57144

$ws.Range("H99").Value = 2437.5
$ws.Range("I99").Value = 2437.5
$ws.Range("K99").Value = 2437.5
$ws.Range("M99").Value = -939.5

$ws.Range("H134").Value = 1358
$ws.Range("J134").Value = 995
$ws.Range("L134").Value = 2985
$ws.Range("N134").Value = -8055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 22000
$ws.Range("I111").Value = 22000
$ws.Range("K111").Value = 22000
$ws.Range("M111").Value = -17910

$ws.Range("H134").Value = 1446.3334
$ws.Range("I134").Value = 1446.3334
$ws.Range("K134").Value = 4339.0002
$ws.Range("M134").Value = -1804.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 56.2
$ws.Range("I7").Value = 56.2
$ws.Range("K7").Value = 168.6
$ws.Range("M7").Value = -56.60000000000002

$ws.Range("H9").Value = 900
$ws.Range("J9").Value = 1000
$ws.Range("L9").Value = 3000
$ws.Range("N9").Value = -3448

$ws.Range("H26").Value = 130
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 130
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 390
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -966

$ws.Range("H34").Value = 584.9524
$ws.Range("J34").Value = 705.6
$ws.Range("L34").Value = 2116.8
$ws.Range("N34").Value = -2284.8

$ws.Range("H40").Value = 105.75
$ws.Range("I40").Value = 77.28570999999999
$ws.Range("K40").Value = 309.14284
$ws.Range("M40").Value = -240.14284

$ws.Range("H92").Value = 549.875
$ws.Range("I92").Value = 483.33334
$ws.Range("J92").Value = 749.5
$ws.Range("K92").Value = 1450.00002
$ws.Range("L92").Value = 2248.5
$ws.Range("M92").Value = -202.0000199999999
$ws.Range("N92").Value = -4744.5

$ws.Range("H97").Value = 530.6
$ws.Range("I97").Value = 538.25
$ws.Range("K97").Value = 1614.75
$ws.Range("M97").Value = -1118.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 147.1
$ws.Range("I2").Value = 55.25
$ws.Range("J2").Value = 208.33333
$ws.Range("K2").Value = 55.25
$ws.Range("L2").Value = 208.33333
$ws.Range("M2").Value = 57.75
$ws.Range("N2").Value = -434.33333

$ws.Range("H43").Value = 11836
$ws.Range("I43").Value = 13672
$ws.Range("K43").Value = 13672
$ws.Range("M43").Value = -13521

$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30970

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H102").Value = 1043.5
$ws.Range("I102").Value = 845.44446
$ws.Range("J102").Value = 1637.6666
$ws.Range("K102").Value = 845.44446
$ws.Range("L102").Value = 1637.6666
$ws.Range("M102").Value = 776.55554
$ws.Range("N102").Value = -4881.6666

$ws.Range("H113").Value = 2749
$ws.Range("J113").Value = 2749
$ws.Range("L113").Value = 2749
$ws.Range("N113").Value = -7089

$ws.Range("H132").Value = 1861.625
$ws.Range("I132").Value = 1832.1666
$ws.Range("J132").Value = 1950
$ws.Range("K132").Value = 5496.4998
$ws.Range("L132").Value = 5850
$ws.Range("M132").Value = -2966.4998
$ws.Range("N132").Value = -10910

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1488.5
$ws.Range("I7").Value = 1612.8182
$ws.Range("J7").Value = 1032.6666
$ws.Range("K7").Value = 1612.8182
$ws.Range("L7").Value = 1032.6666
$ws.Range("M7").Value = -1500.8182
$ws.Range("N7").Value = -1256.6666

$ws.Range("H62").Value = 20249
$ws.Range("J62").Value = 20249
$ws.Range("L62").Value = 20249
$ws.Range("N62").Value = -21497

$ws.Range("H65").Value = 20249
$ws.Range("J65").Value = 20249
$ws.Range("L65").Value = 60747
$ws.Range("N65").Value = -66987

$ws.Range("H93").Value = 1125
$ws.Range("I93").Value = 1125
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1125
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 123
$ws.Range("N93").ClearContents()

$ws.Range("H110").Value = 28249.25
$ws.Range("J110").Value = 28249.25
$ws.Range("L110").Value = 28249.25
$ws.Range("N110").Value = -36429.25

$ws.Range("H126").Value = 1488.5
$ws.Range("I126").Value = 1612.8182
$ws.Range("J126").Value = 1032.6666
$ws.Range("K126").Value = 4838.4546
$ws.Range("L126").Value = 3097.9998
$ws.Range("M126").Value = -2368.4546
$ws.Range("N126").Value = -8037.9998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 29999
$ws.Range("J97").Value = 29999
$ws.Range("L97").Value = 29999
$ws.Range("N97").Value = -31981
